# Update "Economía nacional 2003 a 2021 - Trimestral.xlsx" (Sheet1)
#  - revise several values in the existing last data row (row 74, "01-01-2021")
#  - append a brand-new quarterly row (row 75, "01-04-2021")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Revised figures for row 74 ("01-01-2021")
# ---------------------------------------------------------------------------
$ws.Range("B74").Value  = 16345
$ws.Range("D74").Value  = -5267
$ws.Range("E74").Value  = -956
$ws.Range("F74").Value  = -4311
$ws.Range("G74").Value  = 3610
$ws.Range("I74").Value  = 1494
$ws.Range("K74").Value  = 3506
$ws.Range("M74").Value  = 2612
$ws.Range("N74").Value  = 9463
$ws.Range("O74").Value  = 4679
$ws.Range("Q74").Value  = 4816
$ws.Range("V74").Value  = 6470
$ws.Range("W74").Value  = -1332
$ws.Range("X74").Value  = 17677
$ws.Range("AB74").Value = 3280
$ws.Range("AD74").Value = 1332
$ws.Range("AE74").Value = 4523
$ws.Range("AF74").Value = 1607
$ws.Range("AG74").Value = 2917
$ws.Range("AH74").Value = 7126
$ws.Range("AI74").Value = 3099
$ws.Range("AK74").Value = 4069
$ws.Range("AP74").Value = 7495

# ---------------------------------------------------------------------------
# 2) New row 75 ("01-04-2021")
# ---------------------------------------------------------------------------
# Column A holds a text label that happens to look like a date
# ("01-04-2021"); force the cell to Text first so it is stored as a literal
# shared string instead of being auto-converted to a date serial number,
# then drop the cell format back to the workbook default.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").Style = "Normal"

$ws.Range("B75").Value  = 26655
$ws.Range("C75").Value  = -61
$ws.Range("D75").Value  = 14846
$ws.Range("E75").Value  = 12200
$ws.Range("F75").Value  = 2647
$ws.Range("G75").Value  = 8711
$ws.Range("H75").Value  = 7440
$ws.Range("I75").Value  = 1271
$ws.Range("J75").Value  = 53
$ws.Range("K75").Value  = 10839
$ws.Range("L75").Value  = 7142
$ws.Range("M75").Value  = 3698
$ws.Range("N75").Value  = -5102
$ws.Range("O75").Value  = -1055
$ws.Range("P75").Value  = -732
$ws.Range("Q75").Value  = -3315
$ws.Range("R75").Value  = -8424
$ws.Range("S75").Value  = -560
$ws.Range("T75").Value  = -520
$ws.Range("U75").Value  = -40
$ws.Range("V75").Value  = 6352
$ws.Range("W75").Value  = -1895
$ws.Range("X75").Value  = 28550
$ws.Range("Y75").Value  = 10931
$ws.Range("Z75").Value  = 10270
$ws.Range("AA75").Value = 661
$ws.Range("AB75").Value = 12598
$ws.Range("AC75").Value = 10045
$ws.Range("AD75").Value = 2553
$ws.Range("AE75").Value = 9288
$ws.Range("AF75").Value = 7442
$ws.Range("AG75").Value = 1846
$ws.Range("AH75").Value = 618
$ws.Range("AI75").Value = 1518
$ws.Range("AJ75").Value = -372
$ws.Range("AK75").Value = -527
$ws.Range("AL75").Value = -8424
$ws.Range("AM75").Value = -560
$ws.Range("AN75").Value = -520
$ws.Range("AO75").Value = -40
$ws.Range("AP75").Value = 4099
